# Updated symbol list on Sun Dec 11 21:28:07 UTC 2022 with GitHub Actions
#
# This script reproduces the diff applied to cryptos.xlsx:
#  - Refreshed "Price" (column D) values for most rows.
#  - Row 41 now shows BKEXToken (previously KickToken) and row 43 now shows
#    KickToken (previously BKEXToken) - i.e. the two rows swapped coin
#    identity - each with their own refreshed price / volume label.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, [string] $Text)
    # Force the cell to keep/receive a Text number format so that
    # numeric-looking strings (e.g. "0.006300") are preserved exactly,
    # without Excel silently re-interpreting them as floating point
    # numbers (which would lose trailing zeros / switch to exponential
    # notation for very small values).
    $Range.NumberFormat = "@"
    $Range.Value = $Text
}

# --- Column D (Price) updates -------------------------------------------------
Set-TextValue $ws.Range("D2")  "289.84"
Set-TextValue $ws.Range("D3")  "21.30"
Set-TextValue $ws.Range("D4")  "6.459"
Set-TextValue $ws.Range("D5")  "0.06386"
Set-TextValue $ws.Range("D7")  "1.585"
Set-TextValue $ws.Range("D8")  "6.579"
Set-TextValue $ws.Range("D9")  "0.8279"
Set-TextValue $ws.Range("D11") "0.1687"
Set-TextValue $ws.Range("D12") "0.08712"
Set-TextValue $ws.Range("D13") "0.03659"
Set-TextValue $ws.Range("D14") "0.03205"
Set-TextValue $ws.Range("D15") "0.09199"
Set-TextValue $ws.Range("D16") "3.709"
Set-TextValue $ws.Range("D17") "0.001657"
Set-TextValue $ws.Range("D18") "0.04735"
Set-TextValue $ws.Range("D19") "0.006133"
Set-TextValue $ws.Range("D20") "0.006300"
Set-TextValue $ws.Range("D21") "0.001073"
Set-TextValue $ws.Range("D22") "0.0001602"
Set-TextValue $ws.Range("D23") "3.786"
Set-TextValue $ws.Range("D24") "2.195"
Set-TextValue $ws.Range("D25") "0.3357"
Set-TextValue $ws.Range("D28") "0.0002711"
Set-TextValue $ws.Range("D40") "0.04837"

# --- Row 41: KickToken -> BKEXToken -------------------------------------------
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue $ws.Range("D41") "0.1123"
$ws.Range("E41").Value = "40BKEXTokenBKK"

# --- Row 42 price refresh ------------------------------------------------------
Set-TextValue $ws.Range("D42") "0.003464"

# --- Row 43: BKEXToken -> KickToken -------------------------------------------
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue $ws.Range("D43") "0.007159"
$ws.Range("E43").Value = "42KickTokenKICK"

# --- Remaining column D updates -------------------------------------------------
Set-TextValue $ws.Range("D44") "0.01182"
Set-TextValue $ws.Range("D45") "0.00006913"
Set-TextValue $ws.Range("D46") "0.00000000752"
Set-TextValue $ws.Range("D47") "0.9351"
Set-TextValue $ws.Range("D48") "0.008530"
Set-TextValue $ws.Range("D49") "0.00001904"
Set-TextValue $ws.Range("D50") "0.01243"
